# Auto-generated edit script for SA-HW50.xlsx
# Adds two new simulation methods ("Holden", "Rizzie Spiral"),
# renames "Thomas Hex" -> "Matthies Hex", and reruns the averaging
# simulation so every subsequent row's C:W block shifts down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- extend the label/index columns (A) for the two brand-new rows, copying the existing bold/bordered style ---
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null
$ws.Range("A30").Value = 28
$ws.Range("A31").Value = 29

$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 1.000368325129088
$ws.Range("D4").Value = 0.9994475122671174
$ws.Range("E4").Value = 1.00013812316069
$ws.Range("F4").Value = 0.9998895027303112
$ws.Range("G4").Value = 0.9994475122671174
$ws.Range("H4").Value = 0.9996407712711579
$ws.Range("I4").Value = 1.00013812316069
$ws.Range("J4").Value = 1.000205077294853
$ws.Range("K4").Value = 0.9998812846822938
$ws.Range("L4").Value = 1.000368325129088
$ws.Range("M4").Value = 1.000368325129088
$ws.Range("N4").Value = 1.000368325129088
$ws.Range("O4").Value = 1.00013812316069
$ws.Range("P4").Value = 0.9997928177139036
$ws.Range("Q4").Value = 1.000009703921492
$ws.Range("R4").Value = 0.999984653518965
$ws.Range("S4").Value = 0.999822306703367
$ws.Range("T4").Value = 0.999984653518965
$ws.Range("U4").Value = 0.9999588113097972
$ws.Range("V4").Value = 1.000040714073655
$ws.Range("W4").Value = 0.9999635899620252

$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 1.001519806697356
$ws.Range("D5").Value = 0.9977202967721921
$ws.Range("E5").Value = 1.000569928331517
$ws.Range("F5").Value = 0.999544061183528
$ws.Range("G5").Value = 0.9977202967721921
$ws.Range("H5").Value = 0.9985177267898827
$ws.Range("I5").Value = 1.000569928331517
$ws.Range("J5").Value = 1.000846213706286
$ws.Range("K5").Value = 0.9995101443012864
$ws.Range("L5").Value = 1.001519806697356
$ws.Range("M5").Value = 1.001519806697356
$ws.Range("N5").Value = 1.001519806697356
$ws.Range("O5").Value = 1.000569928331517
$ws.Range("P5").Value = 0.9991451125518543
$ws.Range("Q5").Value = 1.000040036316401
$ws.Range("R5").Value = 0.9999366772670216
$ws.Range("S5").Value = 0.9992667898016651
$ws.Range("T5").Value = 0.9999366772670216
$ws.Range("U5").Value = 0.9998300440255878
$ws.Range("V5").Value = 1.000167996559941
$ws.Range("W5").Value = 0.9998497632641956

$ws.Range("B6").Value = "RotRing OmegaMax-90"
$ws.Range("C6").Value = 1.000592449214823
$ws.Range("D6").Value = 0.9991113251771403
$ws.Range("E6").Value = 1.000222169453412
$ws.Range("F6").Value = 0.9998222657641569
$ws.Range("G6").Value = 0.9991113251771403
$ws.Range("H6").Value = 0.9994221804493494
$ws.Range("I6").Value = 1.000222169453412
$ws.Range("J6").Value = 1.000329867488231
$ws.Range("K6").Value = 0.9998090457806069
$ws.Range("L6").Value = 1.000592449214823
$ws.Range("M6").Value = 1.000592449214823
$ws.Range("N6").Value = 1.000592449214823
$ws.Range("O6").Value = 1.000222169453412
$ws.Range("P6").Value = 0.999666747315276
$ws.Range("Q6").Value = 1.000015607617009
$ws.Range("R6").Value = 0.9999753146151251
$ws.Range("S6").Value = 0.9997141801370529
$ws.Range("T6").Value = 0.9999753146151251
$ws.Range("U6").Value = 0.9999337474064955
$ws.Range("V6").Value = 1.000065487768161
$ws.Range("W6").Value = 0.9999414340976414

$ws.Range("B7").Value = "Equal Angle"
$ws.Range("C7").Value = 1.00047239394092
$ws.Range("D7").Value = 0.9992914085806945
$ws.Range("E7").Value = 1.000177149682999
$ws.Range("F7").Value = 0.9998582821829984
$ws.Range("G7").Value = 0.9992914085806945
$ws.Range("H7").Value = 0.9995392715561964
$ws.Range("I7").Value = 1.000177149682999
$ws.Range("J7").Value = 1.000263021188762
$ws.Range("K7").Value = 0.9998477415994205
$ws.Range("L7").Value = 1.00047239394092
$ws.Range("M7").Value = 1.00047239394092
$ws.Range("N7").Value = 1.00047239394092
$ws.Range("O7").Value = 1.000177149682999
$ws.Range("P7").Value = 0.9997342791318469
$ws.Range("Q7").Value = 1.00001244564121
$ws.Range("R7").Value = 0.9999803174015378
$ws.Range("S7").Value = 0.9997720999543714
$ws.Range("T7").Value = 0.9999803174015378
$ws.Range("U7").Value = 0.9999471734510085
$ws.Range("V7").Value = 1.000052217548991
$ws.Range("W7").Value = 0.9999533023018736

$ws.Range("B8").Value = "Tilt Rotate"
$ws.Range("C8").Value = 1.001529218738082
$ws.Range("D8").Value = 0.9977061614662246
$ws.Range("E8").Value = 1.000573464118299
$ws.Range("F8").Value = 0.9995412344418864
$ws.Range("G8").Value = 0.9977061614662246
$ws.Range("H8").Value = 0.9985085332582859
$ws.Range("I8").Value = 1.000573464118299
$ws.Range("J8").Value = 1.000851441904179
$ws.Range("K8").Value = 0.9995071070939413
$ws.Range("L8").Value = 1.001529218738082
$ws.Range("M8").Value = 1.001529218738082
$ws.Range("N8").Value = 1.001529218738082
$ws.Range("O8").Value = 1.000573464118299
$ws.Range("P8").Value = 0.9991398127922619
$ws.Range("Q8").Value = 1.00004028560612
$ws.Range("R8").Value = 0.9999362814408684
$ws.Range("S8").Value = 0.9992622442261551
$ws.Range("T8").Value = 0.9999362814408684
$ws.Range("U8").Value = 0.9998289878541367
$ws.Range("V8").Value = 1.000169034030926
$ws.Range("W8").Value = 0.9998488281423996

$ws.Range("B9").Value = "CLR"
$ws.Range("C9").Value = 1.000052391116611
$ws.Range("D9").Value = 0.999921413916347
$ws.Range("E9").Value = 1.000019647921593
$ws.Range("F9").Value = 0.9999842831611406
$ws.Range("G9").Value = 0.999921413916347
$ws.Range("H9").Value = 0.9999489050529794
$ws.Range("I9").Value = 1.000019647921593
$ws.Range("J9").Value = 1.00002916962818
$ws.Range("K9").Value = 0.9999831154099968
$ws.Range("L9").Value = 1.000052391116611
$ws.Range("M9").Value = 1.000052391116611
$ws.Range("N9").Value = 1.000052391116611
$ws.Range("O9").Value = 1.000019647921593
$ws.Range("P9").Value = 0.9999705309189701
$ws.Range("Q9").Value = 1.000001381665795
$ws.Range("R9").Value = 0.9999978176515171
$ws.Range("S9").Value = 0.9999747257493125
$ws.Range("T9").Value = 0.9999978176515171
$ws.Range("U9").Value = 0.999994142091137
$ws.Range("V9").Value = 1.000005791896232
$ws.Range("W9").Value = 0.9999948217660553

$ws.Range("B10").Value = "Rizzie Hex"
$ws.Range("C10").Value = 1.00000381009924
$ws.Range("D10").Value = 0.9999942839561042
$ws.Range("E10").Value = 1.000001429719018
$ws.Range("F10").Value = 0.9999988572728581
$ws.Range("G10").Value = 0.9999942839561042
$ws.Range("H10").Value = 0.9999962857943048
$ws.Range("I10").Value = 1.000001429719018
$ws.Range("J10").Value = 1.00000212139633
$ws.Range("K10").Value = 0.9999987734993321
$ws.Range("L10").Value = 1.00000381009924
$ws.Range("M10").Value = 1.00000381009924
$ws.Range("N10").Value = 1.00000381009924
$ws.Range("O10").Value = 1.000001429719018
$ws.Range("P10").Value = 0.9999978568375609
$ws.Range("Q10").Value = 1.000000101609175
$ws.Range("R10").Value = 0.9999998412581205
$ws.Range("S10").Value = 0.9999981623914845
$ws.Range("T10").Value = 0.9999998412581205
$ws.Range("U10").Value = 0.9999995743184233
$ws.Range("V10").Value = 1.000000421474587
$ws.Range("W10").Value = 0.9999996239320255

$ws.Range("B11").Value = "Matthies Hex"
$ws.Range("C11").Value = 1.000084104699112
$ws.Range("D11").Value = 0.9998738430402672
$ws.Range("E11").Value = 1.00003154180094
$ws.Range("F11").Value = 0.999974770014075
$ws.Range("G11").Value = 0.9998738430402672
$ws.Range("H11").Value = 0.9999179742579268
$ws.Range("I11").Value = 1.00003154180094
$ws.Range("J11").Value = 1.000046828368588
$ws.Range("K11").Value = 0.9999728937080901
$ws.Range("L11").Value = 1.000084104699112
$ws.Range("M11").Value = 1.000084104699112
$ws.Range("N11").Value = 1.000084104699112
$ws.Range("O11").Value = 1.00003154180094
$ws.Range("P11").Value = 0.9999526924206038
$ws.Range("Q11").Value = 1.000002217754515
$ws.Range("R11").Value = 0.9999964965134397
$ws.Range("S11").Value = 0.9999594261830992
$ws.Range("T11").Value = 0.9999964965134397
$ws.Range("U11").Value = 0.9999905958121023
$ws.Range("V11").Value = 1.000009297589504
$ws.Range("W11").Value = 0.9999916872112424

$ws.Range("B12").Value = "Tilt Rotate_Partial"
$ws.Range("C12").Value = 1.001551936496884
$ws.Range("D12").Value = 0.9976720846732267
$ws.Range("E12").Value = 1.000581983082369
$ws.Range("F12").Value = 0.9995344188335389
$ws.Range("G12").Value = 0.9976720846732267
$ws.Range("H12").Value = 0.9984863764635588
$ws.Range("I12").Value = 1.000581983082369
$ws.Range("J12").Value = 1.000864090724203
$ws.Range("K12").Value = 0.9994997847790779
$ws.Range("L12").Value = 1.001551936496884
$ws.Range("M12").Value = 1.001551936496884
$ws.Range("N12").Value = 1.001551936496884
$ws.Range("O12").Value = 1.000581983082369
$ws.Range("P12").Value = 0.999127033877798
$ws.Range("Q12").Value = 1.000040883930724
$ws.Range("R12").Value = 0.9999353347508265
$ws.Range("S12").Value = 0.9992512841782246
$ws.Range("T12").Value = 0.9999353347508265
$ws.Range("U12").Value = 0.9998264472578894
$ws.Range("V12").Value = 1.000171545105688
$ws.Range("W12").Value = 0.9998465822669036

$ws.Range("B13").Value = "RotRing OmegaMax-60"
$ws.Range("C13").Value = 1.000301312884267
$ws.Range("D13").Value = 0.9995480303239349
$ws.Range("E13").Value = 1.000112993911815
$ws.Range("F13").Value = 0.9999096074610596
$ws.Range("G13").Value = 0.9995480303239349
$ws.Range("H13").Value = 0.9997061278021158
$ws.Range("I13").Value = 1.000112993911815
$ws.Range("J13").Value = 1.000167767357299
$ws.Range("K13").Value = 0.9999028833104653
$ws.Range("L13").Value = 1.000301312884267
$ws.Range("M13").Value = 1.000301312884267
$ws.Range("N13").Value = 1.000301312884267
$ws.Range("O13").Value = 1.000112993911815
$ws.Range("P13").Value = 0.9998305121178748
$ws.Range("Q13").Value = 1.00000793861114
$ws.Range("R13").Value = 0.9999874457066724
$ws.Range("S13").Value = 0.9998546358487382
$ws.Range("T13").Value = 0.9999874457066724
$ws.Range("U13").Value = 0.9999663051076206
$ws.Range("V13").Value = 1.00003330666295
$ws.Range("W13").Value = 0.9999702146203464

$ws.Range("B14").Value = "Equal Angle_Partial"
$ws.Range("C14").Value = 1.000471754999998
$ws.Range("D14").Value = 0.9992923655894741
$ws.Range("E14").Value = 1.000176910105263
$ws.Range("F14").Value = 0.9998584732421052
$ws.Range("G14").Value = 0.9992923655894741
$ws.Range("H14").Value = 0.9995398946947365
$ws.Range("I14").Value = 1.000176910105263
$ws.Range("J14").Value = 1.000262665168423
$ws.Range("K14").Value = 0.999847947200002
$ws.Range("L14").Value = 1.000471754999998
$ws.Range("M14").Value = 1.000471754999998
$ws.Range("N14").Value = 1.000471754999998
$ws.Range("O14").Value = 1.000176910105263
$ws.Range("P14").Value = 0.9997346378473685
$ws.Range("Q14").Value = 1.000012428652632
$ws.Range("R14").Value = 0.9999803435649115
$ws.Range("S14").Value = 0.9997724076315797
$ws.Range("T14").Value = 0.9999803435649115
$ws.Range("U14").Value = 0.9999472444736841
$ws.Range("V14").Value = 1.000052146578947
$ws.Range("W14").Value = 0.999953365138158

$ws.Range("B15").Value = "Rizzie Hex_Partial"
$ws.Range("C15").Value = 0.9996947719080925
$ws.Range("D15").Value = 1.000457838205996
$ws.Range("E15").Value = 0.9998855396919858
$ws.Range("F15").Value = 1.000091567955177
$ws.Range("G15").Value = 1.000457838205996
$ws.Range("H15").Value = 1.000297694072916
$ws.Range("I15").Value = 0.9998855396919858
$ws.Range("J15").Value = 0.9998300536928426
$ws.Range("K15").Value = 1.00009838075011
$ws.Range("L15").Value = 0.9996947719080925
$ws.Range("M15").Value = 0.9996947719080925
$ws.Range("N15").Value = 0.9996947719080925
$ws.Range("O15").Value = 0.9998855396919858
$ws.Range("P15").Value = 1.000171688948991
$ws.Range("Q15").Value = 0.9999919602210481
$ws.Range("R15").Value = 1.000012716602025
$ws.Range("S15").Value = 1.000147252882698
$ws.Range("T15").Value = 1.000012716602025
$ws.Range("U15").Value = 1.000034132639046
$ws.Range("V15").Value = 0.9999662604928556
$ws.Range("W15").Value = 1.000030173246138

$ws.Range("B16").Value = "ND Single"
$ws.Range("C16").Value = 1.0026813
$ws.Range("D16").Value = 0.9959780300000018
$ws.Range("E16").Value = 1.0010055
$ws.Range("F16").Value = 0.9991956100000006
$ws.Range("G16").Value = 0.9959780300000018
$ws.Range("H16").Value = 0.99738489
$ws.Range("I16").Value = 1.0010055
$ws.Range("J16").Value = 1.001492900000001
$ws.Range("K16").Value = 0.9991357700000003
$ws.Range("L16").Value = 1.0026813
$ws.Range("M16").Value = 1.0026813
$ws.Range("N16").Value = 1.0026813
$ws.Range("O16").Value = 1.0010055
$ws.Range("P16").Value = 0.9984917650000009
$ws.Range("Q16").Value = 1.000070635
$ws.Range("R16").Value = 0.9998882766666674
$ws.Range("S16").Value = 0.998706433333334
$ws.Range("T16").Value = 0.9998882766666674
$ws.Range("U16").Value = 0.9997001500000007
$ws.Range("V16").Value = 1.00029638
$ws.Range("W16").Value = 0.9997349375000004

$ws.Range("B17").Value = "RD Single"
$ws.Range("C17").Value = 1.0043208
$ws.Range("D17").Value = 0.99351878
$ws.Range("E17").Value = 1.0016203
$ws.Range("F17").Value = 0.9987037600000001
$ws.Range("G17").Value = 0.99351878
$ws.Range("H17").Value = 0.99578587
$ws.Range("I17").Value = 1.0016203
$ws.Range("J17").Value = 1.0024058
$ws.Range("K17").Value = 0.99860734
$ws.Range("L17").Value = 1.0043208
$ws.Range("M17").Value = 1.0043208
$ws.Range("N17").Value = 1.0043208
$ws.Range("O17").Value = 1.0016203
$ws.Range("P17").Value = 0.99756954
$ws.Range("Q17").Value = 1.00011382
$ws.Range("R17").Value = 0.99981996
$ws.Range("S17").Value = 0.9979154733333333
$ws.Range("T17").Value = 0.99981996
$ws.Range("U17").Value = 0.999516805
$ws.Range("V17").Value = 1.000477604
$ws.Range("W17").Value = 0.9995728687500001

$ws.Range("B18").Value = "TD Single"
$ws.Range("C18").Value = 1.0041638
$ws.Range("D18").Value = 0.99375433
$ws.Range("E18").Value = 1.0015614
$ws.Range("F18").Value = 0.99875087
$ws.Range("G18").Value = 0.99375433
$ws.Range("H18").Value = 0.9959390299999999
$ws.Range("I18").Value = 1.0015614
$ws.Range("J18").Value = 1.0023183
$ws.Range("K18").Value = 0.99865795
$ws.Range("L18").Value = 1.0041638
$ws.Range("M18").Value = 1.0041638
$ws.Range("N18").Value = 1.0041638
$ws.Range("O18").Value = 1.0015614
$ws.Range("P18").Value = 0.997657865
$ws.Range("Q18").Value = 1.000109675
$ws.Range("R18").Value = 0.99982651
$ws.Range("S18").Value = 0.9979912266666666
$ws.Range("T18").Value = 0.99982651
$ws.Range("U18").Value = 0.99953437
$ws.Range("V18").Value = 1.000460256
$ws.Range("W18").Value = 0.999588385

$ws.Range("B19").Value = "Morris Single"
$ws.Range("C19").Value = 1.0010898
$ws.Range("D19").Value = 0.9983653
$ws.Range("E19").Value = 1.0004087
$ws.Range("F19").Value = 0.9996730599999999
$ws.Range("G19").Value = 0.9983653
$ws.Range("H19").Value = 0.99893711
$ws.Range("I19").Value = 1.0004087
$ws.Range("J19").Value = 1.0006068
$ws.Range("K19").Value = 0.99964874
$ws.Range("L19").Value = 1.0010898
$ws.Range("M19").Value = 1.0010898
$ws.Range("N19").Value = 1.0010898
$ws.Range("O19").Value = 1.0004087
$ws.Range("P19").Value = 0.999387
$ws.Range("Q19").Value = 1.00002872
$ws.Range("R19").Value = 0.9999546
$ws.Range("S19").Value = 0.9994742466666667
$ws.Range("T19").Value = 0.9999546
$ws.Range("U19").Value = 0.999878135
$ws.Range("V19").Value = 1.000120468
$ws.Range("W19").Value = 0.99989227625

$ws.Range("B20").Value = "Ring Perpendicular to ND"
$ws.Range("C20").Value = 1.001050907945205
$ws.Range("D20").Value = 0.9984236389041091
$ws.Range("E20").Value = 1.000394103287671
$ws.Range("F20").Value = 0.9996847317808217
$ws.Range("G20").Value = 0.9984236389041091
$ws.Range("H20").Value = 0.9989750382191785
$ws.Range("I20").Value = 1.000394103287671
$ws.Range("J20").Value = 1.000585128219178
$ws.Range("K20").Value = 0.9996612816438358
$ws.Range("L20").Value = 1.001050907945205
$ws.Range("M20").Value = 1.001050907945205
$ws.Range("N20").Value = 1.001050907945205
$ws.Range("O20").Value = 1.000394103287671
$ws.Range("P20").Value = 0.9994088710958899
$ws.Range("Q20").Value = 1.000027692465753
$ws.Range("R20").Value = 0.9999562167123285
$ws.Range("S20").Value = 0.9994930079452052
$ws.Range("T20").Value = 0.9999562167123285
$ws.Range("U20").Value = 0.9998824829452053
$ws.Range("V20").Value = 1.000116167945205
$ws.Range("W20").Value = 0.9998961166609588

$ws.Range("B21").Value = "Ring Perpendicular to RD"
$ws.Range("C21").Value = 1.001715162631579
$ws.Range("D21").Value = 0.9974272557894737
$ws.Range("E21").Value = 1.000643181578947
$ws.Range("F21").Value = 0.9994854515789473
$ws.Range("G21").Value = 0.9974272557894737
$ws.Range("H21").Value = 0.9983271899999999
$ws.Range("I21").Value = 1.000643181578947
$ws.Range("J21").Value = 1.000954982105263
$ws.Range("K21").Value = 0.99944718
$ws.Range("L21").Value = 1.001715162631579
$ws.Range("M21").Value = 1.001715162631579
$ws.Range("N21").Value = 1.001715162631579
$ws.Range("O21").Value = 1.000643181578947
$ws.Range("P21").Value = 0.9990352186842105
$ws.Range("Q21").Value = 1.000045180789474
$ws.Range("R21").Value = 0.9999285333333333
$ws.Range("S21").Value = 0.999172539122807
$ws.Range("T21").Value = 0.9999285333333333
$ws.Range("U21").Value = 0.999808195
$ws.Range("V21").Value = 1.000189588526316
$ws.Range("W21").Value = 0.9998304481578947

$ws.Range("B22").Value = "Ring Perpendicular to TD"
$ws.Range("C22").Value = 1.001663507894737
$ws.Range("D22").Value = 0.9975047431578947
$ws.Range("E22").Value = 1.000623816842105
$ws.Range("F22").Value = 0.9995009505263158
$ws.Range("G22").Value = 0.9975047431578947
$ws.Range("H22").Value = 0.9983775721052629
$ws.Range("I22").Value = 1.000623816842105
$ws.Range("J22").Value = 1.00092622368421
$ws.Range("K22").Value = 0.9994638268421054
$ws.Range("L22").Value = 1.001663507894737
$ws.Range("M22").Value = 1.001663507894737
$ws.Range("N22").Value = 1.001663507894737
$ws.Range("O22").Value = 1.000623816842105
$ws.Range("P22").Value = 0.99906428
$ws.Range("Q22").Value = 1.000043821842105
$ws.Range("R22").Value = 0.9999306892982456
$ws.Range("S22").Value = 0.9991974622807018
$ws.Range("T22").Value = 0.9999306892982456
$ws.Range("U22").Value = 0.9998139736842105
$ws.Range("V22").Value = 1.000183880526316
$ws.Range("W22").Value = 0.999835557236842

$ws.Range("B23").Value = "OffsetFTD"
$ws.Range("C23").Value = 0.9999974633047395
$ws.Range("D23").Value = 1.000003795180401
$ws.Range("E23").Value = 0.9999990494506927
$ws.Range("F23").Value = 1.000000760577882
$ws.Range("G23").Value = 1.000003795180401
$ws.Range("H23").Value = 1.000002475634174
$ws.Range("I23").Value = 0.9999990494506927
$ws.Range("J23").Value = 0.9999985817290312
$ws.Range("K23").Value = 1.000000829405617
$ws.Range("L23").Value = 0.9999974633047395
$ws.Range("M23").Value = 0.9999974633047395
$ws.Range("N23").Value = 0.9999974633047395
$ws.Range("O23").Value = 0.9999990494506927
$ws.Range("P23").Value = 1.000001422315547
$ws.Range("Q23").Value = 0.9999999394281548
$ws.Range("R23").Value = 1.000000102645278
$ws.Range("S23").Value = 1.000001224678903
$ws.Range("T23").Value = 1.000000102645278
$ws.Range("U23").Value = 1.000000284335362
$ws.Range("V23").Value = 0.9999997201292377
$ws.Range("W23").Value = 1.000000250591654

$ws.Range("B24").Value = "OffsetATD"
$ws.Range("C24").Value = 1.000000664119448
$ws.Range("D24").Value = 0.9999990055467399
$ws.Range("E24").Value = 1.000000246945577
$ws.Range("F24").Value = 0.9999997961638373
$ws.Range("G24").Value = 0.9999990055467399
$ws.Range("H24").Value = 0.999999352278692
$ws.Range("I24").Value = 1.000000246945577
$ws.Range("J24").Value = 1.00000036595987
$ws.Range("K24").Value = 0.9999997910152201
$ws.Range("L24").Value = 1.000000664119448
$ws.Range("M24").Value = 1.000000664119448
$ws.Range("N24").Value = 1.000000664119448
$ws.Range("O24").Value = 1.000000246945577
$ws.Range("P24").Value = 0.9999996262461586
$ws.Range("Q24").Value = 1.000000018980399
$ws.Range("R24").Value = 0.9999999722039217
$ws.Range("S24").Value = 0.999999681169179
$ws.Range("T24").Value = 0.9999999722039217
$ws.Range("U24").Value = 0.9999999269067463
$ws.Range("V24").Value = 1.000000074349287
$ws.Range("W24").Value = 0.9999999336218701

$ws.Range("B25").Value = "OffsetF45"
$ws.Range("C25").Value = 0.9993741563426508
$ws.Range("D25").Value = 1.000938762023138
$ws.Range("E25").Value = 0.9997653092945696
$ws.Range("F25").Value = 1.000187750242171
$ws.Range("G25").Value = 1.000938762023138
$ws.Range("H25").Value = 1.000610391009779
$ws.Range("I25").Value = 0.9997653092945696
$ws.Range("J25").Value = 0.9996515416410945
$ws.Range("K25").Value = 1.000201720479163
$ws.Range("L25").Value = 0.9993741563426508
$ws.Range("M25").Value = 0.9993741563426508
$ws.Range("N25").Value = 0.9993741563426508
$ws.Range("O25").Value = 0.9997653092945696
$ws.Range("P25").Value = 1.000352035658854
$ws.Range("Q25").Value = 0.9999835148868663
$ws.Range("R25").Value = 1.000026075886786
$ws.Range("S25").Value = 1.000301930598957
$ws.Range("T25").Value = 1.000026075886786
$ws.Range("U25").Value = 1.00006998703488
$ws.Range("V25").Value = 0.9999308208964344
$ws.Range("W25").Value = 1.000061867540892

$ws.Range("B26").Value = "OffsetA45"
$ws.Range("C26").Value = 0.9998377853168013
$ws.Range("D26").Value = 1.000243323193166
$ws.Range("E26").Value = 0.9999391702538896
$ws.Range("F26").Value = 1.000048665987001
$ws.Range("G26").Value = 1.000243323193166
$ws.Range("H26").Value = 1.000158218555874
$ws.Range("I26").Value = 0.9999391702538896
$ws.Range("J26").Value = 0.9999096834577995
$ws.Range("K26").Value = 1.000052289871308
$ws.Range("L26").Value = 0.9998377853168013
$ws.Range("M26").Value = 0.9998377853168013
$ws.Range("N26").Value = 0.9998377853168013
$ws.Range("O26").Value = 0.9999391702538896
$ws.Range("P26").Value = 1.000091246723528
$ws.Range("Q26").Value = 0.9999957300625986
$ws.Range("R26").Value = 1.000006759587952
$ws.Range("S26").Value = 1.000078261106121
$ws.Range("T26").Value = 1.000006759587952
$ws.Range("U26").Value = 1.000018142158791
$ws.Range("V26").Value = 0.999982070790393
$ws.Range("W26").Value = 1.000016038361216

$ws.Range("B27").Value = "OffsetFRD"
$ws.Range("C27").Value = 0.9998487039272999
$ws.Range("D27").Value = 1.000226936321122
$ws.Range("E27").Value = 0.999943265563492
$ws.Range("F27").Value = 1.000045395924158
$ws.Range("G27").Value = 1.000226936321122
$ws.Range("H27").Value = 1.000147567502872
$ws.Range("I27").Value = 0.999943265563492
$ws.Range("J27").Value = 0.9999157633476052
$ws.Range("K27").Value = 1.000048770510624
$ws.Range("L27").Value = 0.9998487039272999
$ws.Range("M27").Value = 0.9998487039272999
$ws.Range("N27").Value = 0.9998487039272999
$ws.Range("O27").Value = 0.999943265563492
$ws.Range("P27").Value = 1.000085100942307
$ws.Range("Q27").Value = 0.9999960180370582
$ws.Range("R27").Value = 1.000006301937305
$ws.Range("S27").Value = 1.000072990798413
$ws.Range("T27").Value = 1.000006301937305
$ws.Range("U27").Value = 1.000016919080635
$ws.Range("V27").Value = 0.9999832760499677
$ws.Range("W27").Value = 1.000014958582583

$ws.Range("B28").Value = "OffsetARD"
$ws.Range("C28").Value = 0.9999608861430194
$ws.Range("D28").Value = 1.000058680718968
$ws.Range("E28").Value = 0.9999853301774557
$ws.Range("F28").Value = 1.000011732378739
$ws.Range("G28").Value = 1.000058680718968
$ws.Range("H28").Value = 1.000038155266617
$ws.Range("I28").Value = 0.9999853301774557
$ws.Range("J28").Value = 0.9999782164981107
$ws.Range("K28").Value = 1.000012610808613
$ws.Range("L28").Value = 0.9999608861430194
$ws.Range("M28").Value = 0.9999608861430194
$ws.Range("N28").Value = 0.9999608861430194
$ws.Range("O28").Value = 0.9999853301774557
$ws.Range("P28").Value = 1.000022005448212
$ws.Range("Q28").Value = 0.9999989704930345
$ws.Range("R28").Value = 1.000001632346481
$ws.Range("S28").Value = 1.000018873901679
$ws.Range("T28").Value = 1.000001632346481
$ws.Range("U28").Value = 1.000004376962014
$ws.Range("V28").Value = 0.9999956787982152
$ws.Range("W28").Value = 1.000003867771122

$ws.Range("B29").Value = "Gaussian Quadrature"
$ws.Range("C29").Value = 1.00015980851636
$ws.Range("D29").Value = 0.9997602878264014
$ws.Range("E29").Value = 1.00005991703392
$ws.Range("F29").Value = 0.9999520481835525
$ws.Range("G29").Value = 0.9997602878264014
$ws.Range("H29").Value = 0.9998441278300477
$ws.Range("I29").Value = 1.00005991703392
$ws.Range("J29").Value = 1.000088989405785
$ws.Range("K29").Value = 0.9999484716302068
$ws.Range("L29").Value = 1.00015980851636
$ws.Range("M29").Value = 1.00015980851636
$ws.Range("N29").Value = 1.00015980851636
$ws.Range("O29").Value = 1.00005991703392
$ws.Range("P29").Value = 0.9999101024301609
$ws.Range("Q29").Value = 1.000004194332064
$ws.Range("R29").Value = 0.9999933377922273
$ws.Range("S29").Value = 0.9999228921635095
$ws.Range("T29").Value = 0.9999933377922273
$ws.Range("U29").Value = 0.9999821212517221
$ws.Range("V29").Value = 1.00001765870465
$ws.Range("W29").Value = 0.9999841959325243

$ws.Range("B30").Value = "Michael-CCHex"
$ws.Range("C30").Value = 1.000172719394099
$ws.Range("D30").Value = 0.9997409072362732
$ws.Range("E30").Value = 1.000064773894334
$ws.Range("F30").Value = 0.9999481766818127
$ws.Range("G30").Value = 0.9997409072362732
$ws.Range("H30").Value = 0.9998315450178668
$ws.Range("I30").Value = 1.000064773894334
$ws.Range("J30").Value = 1.000096170020354
$ws.Range("K30").Value = 0.9999443308038072
$ws.Range("L30").Value = 1.000172719394099
$ws.Range("M30").Value = 1.000172719394099
$ws.Range("N30").Value = 1.000172719394099
$ws.Range("O30").Value = 1.000064773894334
$ws.Range("P30").Value = 0.9999028405653037
$ws.Range("Q30").Value = 1.000004552349071
$ws.Range("R30").Value = 0.9999928001749021
$ws.Range("S30").Value = 0.9999166706448048
$ws.Range("T30").Value = 0.9999928001749021
$ws.Range("U30").Value = 0.9999806828321284
$ws.Range("V30").Value = 1.000019090144523
$ws.Range("W30").Value = 0.9999829246178602

$ws.Range("B31").Value = "Michael-SNHex"
$ws.Range("C31").Value = 0.9996825317701024
$ws.Range("D31").Value = 1.000476182415749
$ws.Range("E31").Value = 0.9998809563147099
$ws.Range("F31").Value = 1.000095227772893
$ws.Range("G31").Value = 1.000476182415749
$ws.Range("H31").Value = 1.000309625580092
$ws.Range("I31").Value = 0.9998809563147099
$ws.Range("J31").Value = 0.999823245268298
$ws.Range("K31").Value = 1.000102319457922
$ws.Range("L31").Value = 0.9996825317701024
$ws.Range("M31").Value = 0.9996825317701024
$ws.Range("N31").Value = 0.9996825317701024
$ws.Range("O31").Value = 0.9998809563147099
$ws.Range("P31").Value = 1.00017856936523
$ws.Range("Q31").Value = 0.9999916378863158
$ws.Range("R31").Value = 1.000013223500187
$ws.Range("S31").Value = 1.00015315272946
$ws.Range("T31").Value = 1.000013223500187
$ws.Range("U31").Value = 1.000035497489621
$ws.Range("V31").Value = 0.9999649043457172
$ws.Range("W31").Value = 1.00003138061181

